$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 10: repurpose "Desired annual production" (TEA) into the new
#     "Feedstock capacity" (Feedstock) parameter row ---
$ws.Range("A10").Value = "Feedstock capacity"
$ws.Range("B10").Value = "Feedstock"
$ws.Range("D10").Value = "kg/h"
$ws.Range("E10").Value = 56972
$ws.Range("G10").Formula = "=0.8*H10"
$ws.Range("H10").Value = 56972
$ws.Range("I10").Formula = "=1.2*H10"
$ws.Range("K10").Value = "feedstock.F_mass = x"

# --- Remove the Fermentation TAL yield / titer / productivity rows
#     (old rows 19-21); everything below shifts up and formulas re-point
#     automatically ---
$ws.Rows("19:21").Delete()

# --- Restore the on-screen selection to match the post-edit state ---
[void]$ws.Range("A19:XFD21").Select()
